$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 4 (shifting old rows 4-9 down to 5-10), for the new
# "Session 6" entry (Matt Bombyk - Linking Administrative Data: The IPUMS Experience)
$ws.Rows.Item(4).Insert()

# New row 4 contents
$ws.Range("A4").Value = "1:15PM"
$ws.Range("B4").Value = "1:35PM"
$ws.Range("C4").Value = "20 minutes"
$ws.Range("D4").Value = "Session 6"
$ws.Range("E4").Value = "Linking Administrative Data: The IPUMS Experience"
$ws.Range("F4").Value = "[Matt Bombyk](https://dataifa.github.io/difa-project/comingsoon.html)"

# D4/F4 pick up the time number format from the inserted row; apply explicitly
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("F4").NumberFormat = "h:mm"

# Match the row height used by the other multi-line "session" rows
$ws.Rows.Item(4).RowHeight = 45

# Row 5 (was old row 4) - rename session label and shift start/end times
$ws.Range("A5").Value = "1:35PM"
$ws.Range("B5").Value = "2:35PM"
$ws.Range("D5").Value = "Session 7 (Research Presentations)"

# Row 6 (was old row 5, the Break row) - shift start/end times
$ws.Range("A6").Value = "2:35PM"
$ws.Range("B6").Value = "2:50PM"

# Row 7 (was old row 6, Activity 1) - shift start/end times
$ws.Range("A7").Value = "2:50PM"
$ws.Range("B7").Value = "3:30PM"

# Row 8 (was old row 7, Activity 2) - shift start/end times
$ws.Range("A8").Value = "3:30PM"
$ws.Range("B8").Value = "4:10PM"

# Row 9 (was old row 8, Activity 3) - shift start/end times
$ws.Range("A9").Value = "4:10PM"
$ws.Range("B9").Value = "4:50PM"

# Row 10 (was old row 9, Day 2 Wrap-up) - shift start time
$ws.Range("A10").Value = "4:50PM"

# Update selection to match the new active cell
$ws.Range("D10").Select()
